$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the existing header cell (G1) onto the new header cell (H1),
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values for rows 2-9.
$values = @(0, 1, 0, 0, 1, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
